# Update "想去人数" (interest count) values per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 246
$ws.Range("F3").Value = 231
$ws.Range("F5").Value = 2849
$ws.Range("F6").Value = 62
$ws.Range("F8").Value = 2197
$ws.Range("F9").Value = 317
$ws.Range("G10").Value = 120
$ws.Range("F12").Value = 78
$ws.Range("F13").Value = 2532
$ws.Range("F15").Value = 1312
$ws.Range("F16").Value = 4637
$ws.Range("F18").Value = 4962
$ws.Range("F19").Value = 1528
$ws.Range("F20").Value = 2835
$ws.Range("F21").Value = 3234
$ws.Range("F22").Value = 152
$ws.Range("F23").Value = 1525
$ws.Range("F24").Value = 245
$ws.Range("F25").Value = 833
$ws.Range("F26").Value = 98
$ws.Range("F27").Value = 277
$ws.Range("F28").Value = 945
$ws.Range("F29").Value = 1740
$ws.Range("F30").Value = 113
$ws.Range("F31").Value = 266
$ws.Range("F32").Value = 661
$ws.Range("F33").Value = 154
$ws.Range("F34").Value = 318
$ws.Range("F35").Value = 393

$ws = $wb.Worksheets.Item(2)
$ws.Range("F3").Value = 97
$ws.Range("F8").Value = 87
$ws.Range("F13").Value = 26
$ws.Range("F14").Value = 42

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 97
$ws.Range("F7").Value = 246
$ws.Range("F8").Value = 230
$ws.Range("F11").Value = 2849
$ws.Range("F12").Value = 62
$ws.Range("F13").Value = 2197
$ws.Range("F14").Value = 317
$ws.Range("F15").Value = 87
$ws.Range("G17").Value = 120
$ws.Range("F19").Value = 78
$ws.Range("F21").Value = 2532
$ws.Range("F22").Value = 1313
$ws.Range("F25").Value = 26
$ws.Range("F26").Value = 4637
$ws.Range("F28").Value = 4962
$ws.Range("F29").Value = 1528
$ws.Range("F30").Value = 2835
$ws.Range("F31").Value = 3234
$ws.Range("F32").Value = 152
$ws.Range("F33").Value = 42
$ws.Range("F35").Value = 1525
$ws.Range("F37").Value = 245
$ws.Range("F38").Value = 833
$ws.Range("F39").Value = 98
$ws.Range("F40").Value = 277
$ws.Range("F41").Value = 945
$ws.Range("F43").Value = 1740
$ws.Range("F44").Value = 113
$ws.Range("F45").Value = 266
$ws.Range("F46").Value = 661
$ws.Range("F47").Value = 154
$ws.Range("F48").Value = 318
$ws.Range("F49").Value = 393
